$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1 used to show "ngaytao"; it now shows "tttk", and the old H1 ("tttk") column is removed
$ws.Range("G1").Value = "tttk"
$ws.Range("H1").ClearContents()

# Row 2: huynv
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "huynv"
$ws.Range("C2").Value = "Lương Hoàng Huy"
$ws.Range("D2").Value = "Ho Chi Minh city"
$ws.Range("E2").Value = "huy@gmail.com"
$ws.Range("F2").Value = 10212139
$ws.Range("H2").Value = $true

# Row 3: Huỳnh Tấn Duy (was row with A3=7, now becomes A3=4)
$ws.Range("A3").Value = 4
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Huỳnh Tấn Duy"
$ws.Range("D3").Value = "Ho Chi Minh city"
$ws.Range("E3").Value = "duy@gmail.com"
$ws.Range("F3").Value = 1231355
$ws.Range("H3").Value = $true

# Row 4 (new): Hồ Hữu Đại
$ws.Range("A4").Value = 5
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "Hồ Hữu Đại"
$ws.Range("D4").Value = "Hồ Chí Minh city"
$ws.Range("E4").Value = "dai@gmail.com"
$ws.Range("F4").Value = 11231313
$ws.Range("H4").Value = $true
